# PHOENIX-6395: completed bifurcation of property
# Adds a new "bifurcationDetails" worksheet (mirrors the other *Details
# "data dictionary" sheets used by the functional-test data file) and
# touches a handful of other sheets that were apparently scrolled/clicked
# through while the author was working (expanded used-range, a couple of
# row-height / column-width tweaks).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New sheet: bifurcationDetails (placed after demolitionDetails, i.e.
#    at the very end of the tab strip) with the bifurcation data rows.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bifSheet = $wb.Worksheets.Add($null, $lastSheet)
$bifSheet.Name = "bifurcationDetails"

$bifSheet.Range("A1").Value = "dataName"
$bifSheet.Range("B1").Value = "reasonForCreation"
$bifSheet.Range("C1").Value = "parentAssessmentNo"
$bifSheet.Range("D1").Value = "extentOfSite"
$bifSheet.Range("E1").Value = "occupancyCertificateNumber"

$bifSheet.Range("A2").Value = "bifurcationProperty"
$bifSheet.Range("B2").Value = "BIFURCATION"
$bifSheet.Range("C2").Value = 1016094473
$bifSheet.Range("D2").Value = 5000
$bifSheet.Range("E2").Value = 111

$bifSheet.Columns.Item(1).ColumnWidth = 20.1461
$bifSheet.Columns.Item(2).ColumnWidth = 20.8349
$bifSheet.Columns.Item(3).ColumnWidth = 31.9524
$bifSheet.Columns.Item(4).ColumnWidth = 15.0034
$bifSheet.Columns.Item(5).ColumnWidth = 34.7330
$bifSheet.Columns.Item(6).ColumnWidth = 22.0901

[void]$bifSheet.Range("C10").Select()

# ---------------------------------------------------------------------
# 2. editAssessmentDetails: the used range grew down to row 12.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("editAssessmentDetails")
$ws.Activate()
$ws.Cells.Item(12, 3).Value = "x"
$ws.Cells.Item(12, 3).ClearContents()
[void]$ws.Range("C12").Select()

# ---------------------------------------------------------------------
# 3. editFloorDetails: just re-visited (view scrolled back to A1).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("editFloorDetails")
$ws.Activate()
[void]$ws.Range("D11").Select()

# ---------------------------------------------------------------------
# 4. dataFromWeb: the used range grew down to row 14.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("dataFromWeb")
$ws.Activate()
$ws.Cells.Item(14, 2).Value = "x"
$ws.Cells.Item(14, 2).ClearContents()
[void]$ws.Range("A2").Select()

# ---------------------------------------------------------------------
# 5. registrationDetails: the used range grew down to row 15.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("registrationDetails")
$ws.Activate()
$ws.Cells.Item(15, 17).Value = "x"
$ws.Cells.Item(15, 17).ClearContents()
[void]$ws.Range("F15").Select()

# ---------------------------------------------------------------------
# 6. searchDetails: row 7 height tweaked from 14.65 to 12.8.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("searchDetails")
$ws.Activate()
$ws.Rows.Item(7).RowHeight = 12.8
[void]$ws.Range("B7").Select()

# ---------------------------------------------------------------------
# 7. hearingDetails: first four columns got explicit widths.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("hearingDetails")
$ws.Activate()
$ws.Columns.Item(1).ColumnWidth = 13.0544
$ws.Columns.Item(2).ColumnWidth = 15.1412
$ws.Columns.Item(3).ColumnWidth = 11.9473
$ws.Columns.Item(4).ColumnWidth = 11.5340
[void]$ws.Range("D5").Select()

# ---------------------------------------------------------------------
# 8. assessmentDetails: re-visited, selection left on A1.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("assessmentDetails")
$ws.Activate()
[void]$ws.Range("A1").Select()

# ---------------------------------------------------------------------
# Finally, make the new sheet the active / selected tab (it is the last
# sheet, so this also clears tabSelected on demolitionDetails).
# ---------------------------------------------------------------------
$bifSheet.Activate()
[void]$bifSheet.Range("C10").Select()
